$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1 / sheet1.xml) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(6, 6).Value = 9663
$ws1.Cells.Item(10, 6).Value = 273
$ws1.Cells.Item(12, 6).Value = 476
$ws1.Cells.Item(15, 6).Value = 1208
$ws1.Cells.Item(17, 6).Value = 2997
$ws1.Cells.Item(18, 6).Value = 2267
$ws1.Cells.Item(20, 6).Value = 1967
$ws1.Cells.Item(24, 6).Value = 1565
$ws1.Cells.Item(25, 6).Value = 311
$ws1.Cells.Item(26, 6).Value = 22
$ws1.Cells.Item(27, 6).Value = 190
$ws1.Cells.Item(28, 6).Value = 222
$ws1.Cells.Item(32, 6).Value = 320
$ws1.Cells.Item(35, 6).Value = 154
$ws1.Cells.Item(36, 6).Value = 1546
$ws1.Cells.Item(37, 6).Value = 186
$ws1.Cells.Item(38, 6).Value = 1521
$ws1.Cells.Item(39, 6).Value = 47
$ws1.Cells.Item(40, 6).Value = 359
$ws1.Cells.Item(41, 6).Value = 30
$ws1.Cells.Item(43, 6).Value = 787
$ws1.Cells.Item(44, 6).Value = 78
$ws1.Cells.Item(45, 6).Value = 320

# Sheet "全部类型" (index 4 / sheet4.xml) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(5, 6).Value = 9663
$ws4.Cells.Item(11, 6).Value = 273
$ws4.Cells.Item(13, 6).Value = 476
$ws4.Cells.Item(15, 6).Value = 1208
$ws4.Cells.Item(17, 6).Value = 2997
$ws4.Cells.Item(18, 6).Value = 2267
$ws4.Cells.Item(19, 6).Value = 1967
$ws4.Cells.Item(22, 6).Value = 1565
$ws4.Cells.Item(23, 6).Value = 311
$ws4.Cells.Item(24, 6).Value = 22
$ws4.Cells.Item(25, 6).Value = 190
$ws4.Cells.Item(26, 6).Value = 222
$ws4.Cells.Item(30, 6).Value = 320
$ws4.Cells.Item(36, 6).Value = 154
$ws4.Cells.Item(37, 6).Value = 1546
$ws4.Cells.Item(39, 6).Value = 186
$ws4.Cells.Item(40, 6).Value = 1521
$ws4.Cells.Item(41, 6).Value = 47
$ws4.Cells.Item(43, 6).Value = 359
$ws4.Cells.Item(44, 6).Value = 30
$ws4.Cells.Item(46, 6).Value = 787
$ws4.Cells.Item(47, 6).Value = 78
$ws4.Cells.Item(48, 6).Value = 320
